$d = $word.ActiveDocument

# 1) Title font size: 40 -> 32 (half-points), keep bold, keep paragraph formatting untouched
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Replacement.Font.Size = 16
$find.Execute("MARCO DAVID TOLEDO CANNA", $true, $false, $false, $false, $false, `
              $true, 1, $false, "MARCO DAVID TOLEDO CANNA", 2)

# 2) Update the Front-End experience description text
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Replacement.ClearFormatting()
$oldText = "Form" + [char]0x00E9 + " parte del equipo de desarrollo Front-End, contribuyendo a la creaci" + [char]0x00F3 + "n de interfaces web y m" + [char]0x00F3 + "viles modernas, optimizadas para rendimiento y usabilidad, asegurando buenas pr" + [char]0x00E1 + "cticas de dise" + [char]0x00F1 + "o y desarrollo."
$newText = "Form" + [char]0x00E9 + " parte del equipo de desarrollo Front-End, contribuyendo a la creaci" + [char]0x00F3 + "n de interfaces web y m" + [char]0x00F3 + "viles optimizadas, garantizando un rendimiento eficiente y una experiencia de usuario fluida."
$find2.Execute($oldText, $true, $false, $false, $false, $false, `
               $true, 1, $false, $newText, 2)

# 3) Replace the five empty paragraphs (separators) with a paragraph containing
#    a run of 80 underscores, from the bottom up so indices stay valid.
$separator = "________________________________________________________________________________"
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "" -or $p.Range.Text -eq [char]13 -or $p.Range.Text -eq ([char]13 + [char]7)) {
        $r = $p.Range
        $r.Collapse(1)
        $r.InsertAfter($separator)
    }
}

# 4) Change the "Normal" style font from Calibri to Georgia
$style = $d.Styles.Item("Normal")
$style.Font.Name = "Georgia"
